# edit.ps1 - apply "Update CV with underlined keywords" changes
#
# Helper: underline a substring of a paragraph's text, locating it with
# IndexOf so we never have to hand-count character offsets.
function Underline-Substring {
    param($tr, $paraStart, $fullText, $substring)
    $idx = $fullText.IndexOf($substring)
    if ($idx -lt 0) {
        Write-Host ("Underline-Substring: NOT FOUND [" + $substring + "] in [" + $fullText + "]")
        return
    }
    $rng = $tr.Characters($paraStart + $idx, $substring.Length)
    $rng.Font.Underline = $true
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Shape position / size moves (vertical shift of several boxes that
#    make room for the now-taller experience bullet list).
# ---------------------------------------------------------------------

# Group 1026 (Career Timeline ribbon group)
$grp1026 = $s.Shapes.Item(30)
$grp1026.Top = 306.0750787401575

# TextBox 1056
$tb1056 = $s.Shapes.Item(31)
$tb1056.Top = 648.0383070866142

# TextBox 1057
$tb1057 = $s.Shapes.Item(32)
$tb1057.Top = 516.6821653543307

# TextBox 1058 - also grows taller to fit the rewritten bullets
$tb1058 = $s.Shapes.Item(33)

# TextBox 19
$tb19 = $s.Shapes.Item(40)
$tb19.Top = 457.9906692913386

# ---------------------------------------------------------------------
# 2) TextBox 1058 (Solidity Engineer @ Subspace Labs bullets) - bullets
#    are rewritten & reordered, with key phrases underlined.
# ---------------------------------------------------------------------

$tf1058 = $tb1058.TextFrame
$tr1058 = $tf1058.TextRange

# New plain text for paragraphs 2..8 (COM 1-indexed; paragraph 1 is the
# bold role/date header line, untouched) in their NEW order.
$new2 = "Load testing of Subspace's EVM domain in Async Rust & different Solidity contracts having low to high gas-consuming functions for sending transactions on-chain with custom load value; created Solidity repository monorepo for community."
$new3 = "Built a native cross-chain bridge based on custom LZ contracts between Subspace's EVM and other EVM chains using Solidity & TypeScript; also, extensive Foundry unit & integration testing done."
$new4 = "Did some work on ZK-based Identity solution for native product using TypeScript & Solidity."
$new5 = "ETA Prediction of probabilistic farming reward for a Subspace farmer; released custom GTK4 FE component for Space Acres in Rust."
$new6 = "Plagiarism detection on text embeddings via LSH random projection with Python & Rust."
$new7 = "Worked on SDK development featuring PKI-based identity for Autonomys."
$new8 = "Built a community telegram bot using Rust."

$tr1058.Paragraphs(2,1).Text = $new2
$tr1058.Paragraphs(3,1).Text = $new3
$tr1058.Paragraphs(4,1).Text = $new4
$tr1058.Paragraphs(5,1).Text = $new5
$tr1058.Paragraphs(6,1).Text = $new6
$tr1058.Paragraphs(7,1).Text = $new7
$tr1058.Paragraphs(8,1).Text = $new8

# Apply underline to the key phrases in each rewritten bullet.
Underline-Substring $tr1058 ($tr1058.Paragraphs(2,1).Start) $new2 "Load testing of Subspace's EVM domain"
Underline-Substring $tr1058 ($tr1058.Paragraphs(3,1).Start) $new3 "cross-chain bridge"
Underline-Substring $tr1058 ($tr1058.Paragraphs(4,1).Start) $new4 "ZK-based Identity solution"
Underline-Substring $tr1058 ($tr1058.Paragraphs(5,1).Start) $new5 "custom GTK4 FE component"
Underline-Substring $tr1058 ($tr1058.Paragraphs(6,1).Start) $new6 "Plagiarism detection"
Underline-Substring $tr1058 ($tr1058.Paragraphs(7,1).Start) $new7 "PKI-based identity for Autonomys"

# Resize/reposition TextBox 1058 after the text edits, since spAutoFit
# recalculates height on every text change; our explicit values must be
# applied last so they stick in the saved file.
$tb1058.Top = 326.61964566929134
$tb1058.Height = 129.83066929133858

# ---------------------------------------------------------------------
# 3) TextBox 19 (Jun-Sep 2023 | Self bullets) - underline key phrases in
#    the first two bullets (text & order unchanged otherwise).
# ---------------------------------------------------------------------

$tf19 = $tb19.TextFrame
$tr19 = $tf19.TextRange

$txt2_19 = $tr19.Paragraphs(2,1).Text
Underline-Substring $tr19 ($tr19.Paragraphs(2,1).Start) $txt2_19 "pallet for substrate"

$txt3_19 = $tr19.Paragraphs(3,1).Text
Underline-Substring $tr19 ($tr19.Paragraphs(3,1).Start) $txt3_19 "Coursera Instructor"
